$wb = $excel.ActiveWorkbook

# --- Daily sheet ---
$daily = $wb.Worksheets.Item("Daily")
$daily.Range("G2").Value = 2652.08
$daily.Range("H2").Value = 5890.41
$daily.Range("I2").Value = 686.8200000000001
$daily.Range("J2").Value = 666.11
$daily.Range("L2").Value = 666.11

# --- Hourly sheet ---
$hourly = $wb.Worksheets.Item("Hourly")

# Row 9
$hourly.Range("K9").Value = 0.64
$hourly.Range("M9").Value = 0.64

# Row 10
$hourly.Range("H10").Value = 84.20999999999999
$hourly.Range("I10").Value = 377.5
$hourly.Range("K10").Value = 21.31
$hourly.Range("M10").Value = 21.31

# Row 11
$hourly.Range("K11").Value = 55.31
$hourly.Range("M11").Value = 55.31

# Row 12
$hourly.Range("I12").Value = 712.24
$hourly.Range("K12").Value = 84.70999999999999
$hourly.Range("M12").Value = 84.70999999999999

# Row 13
$hourly.Range("I13").Value = 763.85
$hourly.Range("K13").Value = 104.57
$hourly.Range("M13").Value = 104.57

# Row 14
$hourly.Range("I14").Value = 781.1900000000001
$hourly.Range("K14").Value = 112.57
$hourly.Range("M14").Value = 112.57

# Row 15
$hourly.Range("I15").Value = 769.4400000000001
$hourly.Range("J15").Value = 91.06999999999999
$hourly.Range("K15").Value = 107.5
$hourly.Range("M15").Value = 107.5

# Row 16
$hourly.Range("H16").Value = 355.57
$hourly.Range("I16").Value = 725.05
$hourly.Range("J16").Value = 84.31
$hourly.Range("K16").Value = 89.31
$hourly.Range("M16").Value = 89.31

# Row 17
$hourly.Range("H17").Value = 243.69
$hourly.Range("I17").Value = 632.17
$hourly.Range("J17").Value = 71.53
$hourly.Range("K17").Value = 60.92
$hourly.Range("M17").Value = 60.92

# Row 18
$hourly.Range("H18").Value = 108.78
$hourly.Range("I18").Value = 437
$hourly.Range("K18").Value = 27.2
$hourly.Range("M18").Value = 27.2

# Row 19
$hourly.Range("I19").Value = 62.64
$hourly.Range("K19").Value = 2.07
$hourly.Range("M19").Value = 2.07
